# Training Dashboard update: add new progress row as of 04-Nov-2025
# 1. Insert a new row at position 21 (pushing existing rows 21-36 down to 22-37),
#    copying formatting from the row above (row 20) like Excel's default insert.
# 2. Populate the new row 21 with the new training entry.
# 3. Refresh the "TRAINING DATE"/"EXPIRY DATE" for the two re-trained items
#    (rows 12 and 13, unaffected by the insertion since they are above it).
# 4. Refresh "LAST UPDATE" (I) and "PERIOD TO EXPIRE" (H) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Insert the new row, carrying the formatting of the row above it ---
$ws.Rows.Item(21).EntireRow.Insert()
$ws.Range("A20:K20").Copy()
$ws.Range("A21:K21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Fill in the new row 21 data ---
$ws.Cells.Item(21, 1).Value = 19
$ws.Cells.Item(21, 2).Value = "Diagnosis Of Beckoff Module And Troubleshooting Guide (Cargo Trainings)"
$ws.Cells.Item(21, 3).Value = "CARGO"
$ws.Cells.Item(21, 4).Value = "LSME-CRG-M-012"
$ws.Cells.Item(21, 5).Value = "EQUIPMENT MANUAL"
$ws.Cells.Item(21, 6).Value = "21-Oct-2025"
$ws.Cells.Item(21, 7).Value = "21-Oct-2027"
$ws.Cells.Item(21, 8).Value = 715
$ws.Cells.Item(21, 9).Value = "04-Nov-2025"
$ws.Cells.Item(21, 10).Value = "VALID"
$ws.Cells.Item(21, 11).Value = ""

# --- 3. Update training/expiry dates for re-trained items (rows 12 & 13) ---
$ws.Cells.Item(12, 6).Value = "20-Oct-2025"
$ws.Cells.Item(12, 7).Value = "20-Oct-2027"

$ws.Cells.Item(13, 6).Value = "20-Oct-2025"
$ws.Cells.Item(13, 7).Value = "20-Oct-2027"

# --- 4. Refresh LAST UPDATE (I) and PERIOD TO EXPIRE (H) for every data row ---
$updates = @(
    @{ Row=3;  H=703;  I="04-Nov-2025" },
    @{ Row=4;  H=704;  I="04-Nov-2025" },
    @{ Row=5;  H=712;  I="04-Nov-2025" },
    @{ Row=6;  H=702;  I="04-Nov-2025" },
    @{ Row=7;  H=712;  I="04-Nov-2025" },
    @{ Row=8;  H=368;  I="04-Nov-2025" },
    @{ Row=9;  H=704;  I="04-Nov-2025" },
    @{ Row=10; H=712;  I="04-Nov-2025" },
    @{ Row=11; H=703;  I="04-Nov-2025" },
    @{ Row=12; H=714;  I="04-Nov-2025" },
    @{ Row=13; H=714;  I="04-Nov-2025" },
    @{ Row=14; H=361;  I="04-Nov-2025" },
    @{ Row=15; H=362;  I="04-Nov-2025" },
    @{ Row=16; H=705;  I="04-Nov-2025" },
    @{ Row=17; H=425;  I="04-Nov-2025" },
    @{ Row=18; H=424;  I="04-Nov-2025" },
    @{ Row=19; H=423;  I="04-Nov-2025" },
    @{ Row=20; H=424;  I="04-Nov-2025" },
    @{ Row=21; H=715;  I="04-Nov-2025" },
    @{ Row=22; H=-28;  I="04-Nov-2025" },
    @{ Row=23; H=-152; I="04-Nov-2025" },
    @{ Row=24; H=-104; I="04-Nov-2025" },
    @{ Row=25; H=137;  I="04-Nov-2025" },
    @{ Row=26; H=-48;  I="04-Nov-2025" },
    @{ Row=27; H=136;  I="04-Nov-2025" },
    @{ Row=28; H=151;  I="04-Nov-2025" },
    @{ Row=29; H=151;  I="04-Nov-2025" },
    @{ Row=30; H=263;  I="04-Nov-2025" },
    @{ Row=31; H=263;  I="04-Nov-2025" },
    @{ Row=32; H=263;  I="04-Nov-2025" },
    @{ Row=33; H=263;  I="04-Nov-2025" },
    @{ Row=34; H=347;  I="04-Nov-2025" },
    @{ Row=35; H=284;  I="04-Nov-2025" },
    @{ Row=36; H=284;  I="04-Nov-2025" },
    @{ Row=37; H=604;  I="04-Nov-2025" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 8).Value = $u.H
    $ws.Cells.Item($u.Row, 9).Value = $u.I
}
